# Auto-generated script to apply odds updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Cells.Item(3, 13).Value = 1.14
$ws.Cells.Item(3, 15).Value = 1.57
$ws.Cells.Item(3, 19).Value = 2.87
$ws.Cells.Item(3, 20).Value = 1.37
$ws.Cells.Item(3, 24).Value = 1.13

# Row 4
$ws.Cells.Item(4, 7).Value = 2.45
$ws.Cells.Item(4, 9).Value = 3.5
$ws.Cells.Item(4, 10).Value = 3.4
$ws.Cells.Item(4, 13).Value = 1.17
$ws.Cells.Item(4, 15).Value = 1.67
$ws.Cells.Item(4, 17).Value = 2.38
$ws.Cells.Item(4, 18).Value = 1.59
$ws.Cells.Item(4, 20).Value = 1.3
$ws.Cells.Item(4, 21).Value = 5.8
$ws.Cells.Item(4, 22).Value = 1.14
$ws.Cells.Item(4, 24).Value = 1.1
$ws.Cells.Item(4, 27).Value = 2.5
$ws.Cells.Item(4, 28).Value = 1.5
$ws.Cells.Item(4, 30).Value = 9.5
$ws.Cells.Item(4, 32).Value = 23
$ws.Cells.Item(4, 40).Value = 7

# Row 5
$ws.Cells.Item(5, 15).Value = 1.72

# Row 7
$ws.Cells.Item(7, 7).Value = 1.8
$ws.Cells.Item(7, 8).Value = 3.1
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 2.42
$ws.Cells.Item(7, 11).Value = 1.98
$ws.Cells.Item(7, 12).Value = 5.3
$ws.Cells.Item(7, 15).Value = 1.47
$ws.Cells.Item(7, 16).Value = 2.5
$ws.Cells.Item(7, 19).Value = 2.37
$ws.Cells.Item(7, 20).Value = 1.52
$ws.Cells.Item(7, 23).Value = 4.2
$ws.Cells.Item(7, 26).Value = 2.4
$ws.Cells.Item(7, 30).Value = 7.1
$ws.Cells.Item(7, 31).Value = 8.75
$ws.Cells.Item(7, 32).Value = 14
$ws.Cells.Item(7, 33).Value = 17.5
$ws.Cells.Item(7, 34).Value = 40
$ws.Cells.Item(7, 36).Value = 6.2
$ws.Cells.Item(7, 37).Value = 19
$ws.Cells.Item(7, 40).Value = 10.5
$ws.Cells.Item(7, 41).Value = 27
$ws.Cells.Item(7, 42).Value = 16.5
$ws.Cells.Item(7, 43).Value = 100
$ws.Cells.Item(7, 44).Value = 60
$ws.Cells.Item(7, 45).Value = 70

# Row 8
$ws.Cells.Item(8, 8).Value = 2.57
$ws.Cells.Item(8, 9).Value = 2.82
$ws.Cells.Item(8, 10).Value = 3.8
$ws.Cells.Item(8, 12).Value = 3.6
$ws.Cells.Item(8, 13).Value = 1.17
$ws.Cells.Item(8, 14).Value = 4.45
$ws.Cells.Item(8, 15).Value = 1.7
$ws.Cells.Item(8, 16).Value = 2.05
$ws.Cells.Item(8, 19).Value = 3
$ws.Cells.Item(8, 20).Value = 1.34
$ws.Cells.Item(8, 23).Value = 5.5
$ws.Cells.Item(8, 24).Value = 1.11
$ws.Cells.Item(8, 25).Value = 1.65
$ws.Cells.Item(8, 26).Value = 2.1
$ws.Cells.Item(8, 27).Value = 2.27
$ws.Cells.Item(8, 28).Value = 1.57
$ws.Cells.Item(8, 29).Value = 6.1
$ws.Cells.Item(8, 30).Value = 13.5
$ws.Cells.Item(8, 31).Value = 11.75
$ws.Cells.Item(8, 34).Value = 60
$ws.Cells.Item(8, 35).Value = 4.45
$ws.Cells.Item(8, 36).Value = 5.3
$ws.Cells.Item(8, 37).Value = 19.5
$ws.Cells.Item(8, 38).Value = 150
$ws.Cells.Item(8, 40).Value = 5.8
$ws.Cells.Item(8, 41).Value = 12
$ws.Cells.Item(8, 42).Value = 11.5
$ws.Cells.Item(8, 44).Value = 35
$ws.Cells.Item(8, 45).Value = 60

# Row 9
$ws.Cells.Item(9, 7).Value = 1.37
$ws.Cells.Item(9, 8).Value = 4.2
$ws.Cells.Item(9, 9).Value = 9
$ws.Cells.Item(9, 10).Value = 1.88
$ws.Cells.Item(9, 11).Value = 2.25
$ws.Cells.Item(9, 13).Value = 1.05
$ws.Cells.Item(9, 14).Value = 7.8
$ws.Cells.Item(9, 15).Value = 1.25
$ws.Cells.Item(9, 16).Value = 3.55
$ws.Cells.Item(9, 19).Value = 1.75
$ws.Cells.Item(9, 20).Value = 1.95
$ws.Cells.Item(9, 23).Value = 2.8
$ws.Cells.Item(9, 24).Value = 1.38
$ws.Cells.Item(9, 25).Value = 1.38
$ws.Cells.Item(9, 26).Value = 2.82
$ws.Cells.Item(9, 27).Value = 2.02
$ws.Cells.Item(9, 28).Value = 1.72
$ws.Cells.Item(9, 29).Value = 6.2
$ws.Cells.Item(9, 30).Value = 6.1
$ws.Cells.Item(9, 33).Value = 11.5
$ws.Cells.Item(9, 34).Value = 28
$ws.Cells.Item(9, 35).Value = 7.8
$ws.Cells.Item(9, 36).Value = 8.5
$ws.Cells.Item(9, 37).Value = 20
$ws.Cells.Item(9, 38).Value = 100
$ws.Cells.Item(9, 39).Value = 800
$ws.Cells.Item(9, 40).Value = 23
$ws.Cells.Item(9, 41).Value = 65
$ws.Cells.Item(9, 45).Value = 80

# Row 10
$ws.Cells.Item(10, 7).Value = 1.75
$ws.Cells.Item(10, 9).Value = 3.8
$ws.Cells.Item(10, 10).Value = 2.38
$ws.Cells.Item(10, 12).Value = 4.33
$ws.Cells.Item(10, 15).Value = 1.17
$ws.Cells.Item(10, 16).Value = 5
$ws.Cells.Item(10, 19).Value = 1.53
$ws.Cells.Item(10, 20).Value = 2.4
$ws.Cells.Item(10, 21).Value = 1.85
$ws.Cells.Item(10, 22).Value = 1.95
$ws.Cells.Item(10, 23).Value = 2.25
$ws.Cells.Item(10, 24).Value = 1.57
$ws.Cells.Item(10, 27).Value = 1.53
$ws.Cells.Item(10, 28).Value = 2.38
$ws.Cells.Item(10, 29).Value = 10
$ws.Cells.Item(10, 30).Value = 10
$ws.Cells.Item(10, 31).Value = 8.5
$ws.Cells.Item(10, 32).Value = 15
$ws.Cells.Item(10, 34).Value = 21
$ws.Cells.Item(10, 35).Value = 17
$ws.Cells.Item(10, 36).Value = 7.5
$ws.Cells.Item(10, 37).Value = 13
$ws.Cells.Item(10, 38).Value = 41
$ws.Cells.Item(10, 39).Value = 126
$ws.Cells.Item(10, 40).Value = 17
$ws.Cells.Item(10, 41).Value = 23
$ws.Cells.Item(10, 42).Value = 13
$ws.Cells.Item(10, 44).Value = 29
$ws.Cells.Item(10, 45).Value = 29

# Row 12
$ws.Cells.Item(12, 7).Value = 1.2
$ws.Cells.Item(12, 8).Value = 6.25
$ws.Cells.Item(12, 9).Value = 10
$ws.Cells.Item(12, 10).Value = 1.57
$ws.Cells.Item(12, 12).Value = 9
$ws.Cells.Item(12, 19).Value = 1.44
$ws.Cells.Item(12, 20).Value = 2.63
$ws.Cells.Item(12, 27).Value = 2.1
$ws.Cells.Item(12, 28).Value = 1.63
$ws.Cells.Item(12, 31).Value = 10
$ws.Cells.Item(12, 32).Value = 7.5
$ws.Cells.Item(12, 37).Value = 26
$ws.Cells.Item(12, 38).Value = 67
$ws.Cells.Item(12, 42).Value = 29
$ws.Cells.Item(12, 43).Value = 126
$ws.Cells.Item(12, 44).Value = 67
$ws.Cells.Item(12, 45).Value = 51

# Row 13
$ws.Cells.Item(13, 7).Value = 1.67
$ws.Cells.Item(13, 9).Value = 3.9
$ws.Cells.Item(13, 11).Value = 2.38
$ws.Cells.Item(13, 19).Value = 1.57
$ws.Cells.Item(13, 20).Value = 2.35
$ws.Cells.Item(13, 21).Value = 1.95
$ws.Cells.Item(13, 22).Value = 1.85
$ws.Cells.Item(13, 23).Value = 2.38
$ws.Cells.Item(13, 24).Value = 1.53
$ws.Cells.Item(13, 27).Value = 1.58
$ws.Cells.Item(13, 40).Value = 15
$ws.Cells.Item(13, 41).Value = 23
$ws.Cells.Item(13, 42).Value = 13
$ws.Cells.Item(13, 45).Value = 29

# Row 14
$ws.Cells.Item(14, 7).Value = 4.75
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(14, 9).Value = 1.57
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 27).Value = 1.8
$ws.Cells.Item(14, 35).Value = 12
$ws.Cells.Item(14, 37).Value = 17
$ws.Cells.Item(14, 39).Value = 700
$ws.Cells.Item(14, 43).Value = 12

# Row 15
$ws.Cells.Item(15, 13).Value = 1.01
$ws.Cells.Item(15, 14).Value = 19
$ws.Cells.Item(15, 15).Value = 1.11
$ws.Cells.Item(15, 21).Value = 1.83
$ws.Cells.Item(15, 22).Value = 2.03
$ws.Cells.Item(15, 24).Value = 1.58
$ws.Cells.Item(15, 25).Value = 1.25
$ws.Cells.Item(15, 26).Value = 3.75
$ws.Cells.Item(15, 27).Value = 1.5
$ws.Cells.Item(15, 28).Value = 2.5
$ws.Cells.Item(15, 29).Value = 11
$ws.Cells.Item(15, 30).Value = 11
$ws.Cells.Item(15, 35).Value = 19

# Row 16
$ws.Cells.Item(16, 14).Value = 13
$ws.Cells.Item(16, 27).Value = 1.63

# Row 18
$ws.Cells.Item(18, 7).Value = 1.42
$ws.Cells.Item(18, 8).Value = 4.25
$ws.Cells.Item(18, 9).Value = 6.8
$ws.Cells.Item(18, 10).Value = 1.91
$ws.Cells.Item(18, 11).Value = 2.27
$ws.Cells.Item(18, 12).Value = 6.4
$ws.Cells.Item(18, 15).Value = 1.26
$ws.Cells.Item(18, 16).Value = 3.15
$ws.Cells.Item(18, 19).Value = 1.78
$ws.Cells.Item(18, 20).Value = 1.83
$ws.Cells.Item(18, 23).Value = 2.82
$ws.Cells.Item(18, 24).Value = 1.32
$ws.Cells.Item(18, 27).Value = 2.02
$ws.Cells.Item(18, 28).Value = 1.62
$ws.Cells.Item(18, 29).Value = 6.2
$ws.Cells.Item(18, 30).Value = 6.1
$ws.Cells.Item(18, 31).Value = 8.5
$ws.Cells.Item(18, 32).Value = 8.75
$ws.Cells.Item(18, 33).Value = 12
$ws.Cells.Item(18, 34).Value = 32
$ws.Cells.Item(18, 36).Value = 8.5
$ws.Cells.Item(18, 37).Value = 22
$ws.Cells.Item(18, 38).Value = 120
$ws.Cells.Item(18, 39).Value = 900
$ws.Cells.Item(18, 40).Value = 16.5
$ws.Cells.Item(18, 41).Value = 45
$ws.Cells.Item(18, 42).Value = 22
$ws.Cells.Item(18, 43).Value = 150
$ws.Cells.Item(18, 44).Value = 80
$ws.Cells.Item(18, 45).Value = 80
